$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 27.852944
$ws.Range("H2").Value = 83.558832
$ws.Range("I2").Value = 0.2559209115167818
$ws.Range("J2").Value = 0.2559209115167818
$ws.Range("M2").Value = 8.226724333333333
$ws.Range("N2").Value = 24.680173
$ws.Range("O2").Value = 0.06198126651953669
$ws.Range("P2").Value = 0.06198126651953669
$ws.Range("Q2").Value = 229.1384921597706
$ws.Range("R2").Value = 2062.246429437936
$ws.Range("S2").Value = 0.01586230222464442
$ws.Range("T2").Value = 0.01586230222464442
$ws.Range("G3").Value = 27.852944
$ws.Range("H3").Value = 83.558832
$ws.Range("I3").Value = 0.2559209115167818
$ws.Range("J3").Value = 0.2559209115167818
$ws.Range("O3").Value = 0.6623065855236785
$ws.Range("P3").Value = 0.6623065855236785
$ws.Range("Q3").Value = 2448.480659983717
$ws.Range("R3").Value = 22036.32593985345
$ws.Range("S3").Value = 0.1694981050707872
$ws.Range("T3").Value = 0.1694981050707872
$ws.Range("G4").Value = 27.852944
$ws.Range("H4").Value = 83.558832
$ws.Range("I4").Value = 0.2559209115167818
$ws.Range("J4").Value = 0.2559209115167818
$ws.Range("M4").Value = 36.43008433333333
$ws.Range("N4").Value = 109.290253
$ws.Range("O4").Value = 0.2744692388979848
$ws.Range("P4").Value = 0.2744692388979848
$ws.Range("Q4").Value = 1014.685098851611
$ws.Range("R4").Value = 9132.165889664497
$ws.Range("S4").Value = 0.0702424178020896
$ws.Range("T4").Value = 0.0702424178020896
$ws.Range("G5").Value = 27.852944
$ws.Range("H5").Value = 83.558832
$ws.Range("I5").Value = 0.2559209115167818
$ws.Range("J5").Value = 0.2559209115167818
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1649703333333333
$ws.Range("N5").Value = 0.494911
$ws.Range("O5").Value = 0.00124290905879997
$ws.Range("P5").Value = 0.00124290905879997
$ws.Range("Q5").Value = 4.594909455994666
$ws.Range("R5").Value = 41.354185103952
$ws.Range("S5").Value = 0.0003180864192605535
$ws.Range("T5").Value = 0.0003180864192605535
$ws.Range("I6").Value = 0.3112048767201538
$ws.Range("J6").Value = 0.3112048767201538
$ws.Range("M6").Value = 8.226724333333333
$ws.Range("N6").Value = 24.680173
$ws.Range("O6").Value = 0.06198126651953669
$ws.Range("P6").Value = 0.06198126651953669
$ws.Range("Q6").Value = 278.6369264699471
$ws.Range("R6").Value = 2507.732338229524
$ws.Range("S6").Value = 0.01928887240617141
$ws.Range("T6").Value = 0.01928887240617141
$ws.Range("I7").Value = 0.3112048767201538
$ws.Range("J7").Value = 0.3112048767201538
$ws.Range("O7").Value = 0.6623065855236785
$ws.Range("P7").Value = 0.6623065855236785
$ws.Range("R7").Value = 26796.60703313621
$ws.Range("S7").Value = 0.2061130392988424
$ws.Range("T7").Value = 0.2061130392988424
$ws.Range("I8").Value = 0.3112048767201538
$ws.Range("J8").Value = 0.3112048767201538
$ws.Range("M8").Value = 36.43008433333333
$ws.Range("N8").Value = 109.290253
$ws.Range("O8").Value = 0.2744692388979848
$ws.Range("P8").Value = 0.2744692388979848
$ws.Range("Q8").Value = 1233.877095960507
$ws.Range("R8").Value = 11104.89386364457
$ws.Range("S8").Value = 0.08541616565472182
$ws.Range("T8").Value = 0.08541616565472182
$ws.Range("I9").Value = 0.3112048767201538
$ws.Range("J9").Value = 0.3112048767201538
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.1649703333333333
$ws.Range("N9").Value = 0.494911
$ws.Range("O9").Value = 0.00124290905879997
$ws.Range("P9").Value = 0.00124290905879997
$ws.Range("Q9").Value = 5.587500538029778
$ws.Range("R9").Value = 50.287504842268
$ws.Range("S9").Value = 0.000386799360418207
$ws.Range("T9").Value = 0.000386799360418207
$ws.Range("G10").Value = 30.14135433333333
$ws.Range("H10").Value = 90.42406299999999
$ws.Range("I10").Value = 0.2769474880406526
$ws.Range("J10").Value = 0.2769474880406526
$ws.Range("M10").Value = 8.226724333333333
$ws.Range("N10").Value = 24.680173
$ws.Range("O10").Value = 0.06198126651953669
$ws.Range("P10").Value = 0.06198126651953669
$ws.Range("Q10").Value = 247.9646131336554
$ws.Range("R10").Value = 2231.681518202899
$ws.Range("S10").Value = 0.01716555606816389
$ws.Range("T10").Value = 0.01716555606816389
$ws.Range("G11").Value = 30.14135433333333
$ws.Range("H11").Value = 90.42406299999999
$ws.Range("I11").Value = 0.2769474880406526
$ws.Range("J11").Value = 0.2769474880406526
$ws.Range("O11").Value = 0.6623065855236785
$ws.Range("P11").Value = 0.6623065855236785
$ws.Range("Q11").Value = 2649.648925832869
$ws.Range("R11").Value = 23846.84033249583
$ws.Range("S11").Value = 0.1834241451735644
$ws.Range("T11").Value = 0.1834241451735644
$ws.Range("G12").Value = 30.14135433333333
$ws.Range("H12").Value = 90.42406299999999
$ws.Range("I12").Value = 0.2769474880406526
$ws.Range("J12").Value = 0.2769474880406526
$ws.Range("M12").Value = 36.43008433333333
$ws.Range("N12").Value = 109.290253
$ws.Range("O12").Value = 0.2744692388979848
$ws.Range("P12").Value = 0.2744692388979848
$ws.Range("Q12").Value = 1098.052080284215
$ws.Range("R12").Value = 9882.468722557938
$ws.Range("S12").Value = 0.07601356625722668
$ws.Range("T12").Value = 0.07601356625722666
$ws.Range("G13").Value = 30.14135433333333
$ws.Range("H13").Value = 90.42406299999999
$ws.Range("I13").Value = 0.2769474880406526
$ws.Range("J13").Value = 0.2769474880406526
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.1649703333333333
$ws.Range("N13").Value = 0.494911
$ws.Range("O13").Value = 0.00124290905879997
$ws.Range("P13").Value = 0.00124290905879997
$ws.Range("Q13").Value = 4.97242927148811
$ws.Range("R13").Value = 44.75186344339299
$ws.Range("S13").Value = 0.0003442205416976234
$ws.Range("T13").Value = 0.0003442205416976233
$ws.Range("G14").Value = 16.970158
$ws.Range("H14").Value = 50.910474
$ws.Range("I14").Value = 0.1559267237224118
$ws.Range("J14").Value = 0.1559267237224118
$ws.Range("M14").Value = 8.226724333333333
$ws.Range("N14").Value = 24.680173
$ws.Range("O14").Value = 0.06198126651953669
$ws.Range("P14").Value = 0.06198126651953669
$ws.Range("Q14").Value = 139.6088117591113
$ws.Range("R14").Value = 1256.479305832002
$ws.Range("S14").Value = 0.00966453582055697
$ws.Range("T14").Value = 0.009664535820556968
$ws.Range("G15").Value = 16.970158
$ws.Range("H15").Value = 50.910474
$ws.Range("I15").Value = 0.1559267237224118
$ws.Range("J15").Value = 0.1559267237224118
$ws.Range("O15").Value = 0.6623065855236785
$ws.Range("P15").Value = 0.6623065855236785
$ws.Range("Q15").Value = 1491.802936876905
$ws.Range("R15").Value = 13426.22643189214
$ws.Range("S15").Value = 0.1032712959804845
$ws.Range("T15").Value = 0.1032712959804845
$ws.Range("G16").Value = 16.970158
$ws.Range("H16").Value = 50.910474
$ws.Range("I16").Value = 0.1559267237224118
$ws.Range("J16").Value = 0.1559267237224118
$ws.Range("M16").Value = 36.43008433333333
$ws.Range("N16").Value = 109.290253
$ws.Range("O16").Value = 0.2744692388979848
$ws.Range("P16").Value = 0.2744692388979848
$ws.Range("Q16").Value = 618.2242870899914
$ws.Range("R16").Value = 5564.018583809922
$ws.Range("S16").Value = 0.04279708918394672
$ws.Range("T16").Value = 0.04279708918394671
$ws.Range("G17").Value = 16.970158
$ws.Range("H17").Value = 50.910474
$ws.Range("I17").Value = 0.1559267237224118
$ws.Range("J17").Value = 0.1559267237224118
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.1649703333333333
$ws.Range("N17").Value = 0.494911
$ws.Range("O17").Value = 0.00124290905879997
$ws.Range("P17").Value = 0.00124290905879997
$ws.Range("Q17").Value = 2.799572621979333
$ws.Range("R17").Value = 25.196153597814
$ws.Range("S17").Value = 0.0001938027374235857
$ws.Range("T17").Value = 0.0001938027374235857
